$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4-7 (DK1_SmallDecentral_HS, DK2_Central_HS, DK2_LargeDecentral_HS, DK2_SmallDecentral_HS)
$ws.Rows("4:7").Delete()

# Update row 2 values: DK1_Central_HS -> id_DK_Central_HS, DK1 -> DK
$ws.Range("A2").Value = "id_DK_Central_HS"
$ws.Range("C2").Value = "id_DK_Central_HS"
$ws.Range("E2").Value = "id_DK_Central_HS"
$ws.Range("F2").Value = "DK"

# Update row 3 values: DK1_LargeDecentral_HS -> id_DK_Decentral_HS, DK1 -> DK
$ws.Range("A3").Value = "id_DK_Decentral_HS"
$ws.Range("C3").Value = "id_DK_Decentral_HS"
$ws.Range("E3").Value = "id_DK_Decentral_HS"
$ws.Range("F3").Value = "DK"
